# Refresh cached market-board figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# for the leve-profit rows across all eight crafting-job sheets, as produced by
# the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1592.9412
$ws.Range("I15").Value = 1592.9412
$ws.Range("K15").Value = 4778.8236
$ws.Range("M15").Value = -4609.8236
$ws.Range("H21").Value = 17471.588
$ws.Range("I21").Value = 26008.5
$ws.Range("J21").Value = 16333.333
$ws.Range("K21").Value = 26008.5
$ws.Range("L21").Value = 16333.333
$ws.Range("M21").Value = -25540.5
$ws.Range("N21").Value = -17269.333
$ws.Range("H23").Value = 17471.588
$ws.Range("I23").Value = 26008.5
$ws.Range("J23").Value = 16333.333
$ws.Range("K23").Value = 26008.5
$ws.Range("L23").Value = 16333.333
$ws.Range("M23").Value = -25774.5
$ws.Range("N23").Value = -16801.333
$ws.Range("H112").Value = 2853.484
$ws.Range("J112").Value = 3027.862
$ws.Range("L112").Value = 9083.585999999999
$ws.Range("N112").Value = -11299.586
$ws.Range("H132").Value = 2219.3394
$ws.Range("I132").Value = 1628.2821
$ws.Range("J132").Value = 3575.2942
$ws.Range("K132").Value = 4884.846299999999
$ws.Range("L132").Value = 10725.8826
$ws.Range("M132").Value = -2354.846299999999
$ws.Range("N132").Value = -15785.8826
$ws.Range("H135").Value = 732.4595
$ws.Range("I135").Value = 689.36365
$ws.Range("J135").Value = 1088
$ws.Range("K135").Value = 6204.27285
$ws.Range("L135").Value = 9792
$ws.Range("M135").Value = -3669.27285
$ws.Range("N135").Value = -14862
$ws.Range("H137").Value = 1715.7778
$ws.Range("I137").Value = 1156.7222
$ws.Range("K137").Value = 3470.1666
$ws.Range("M137").Value = -920.1665999999996
$ws.Range("H138").Value = 2611.0745
$ws.Range("I138").Value = 1494.9375
$ws.Range("J138").Value = 3775.739
$ws.Range("K138").Value = 4484.8125
$ws.Range("L138").Value = 11327.217
$ws.Range("M138").Value = 655.1875
$ws.Range("N138").Value = -21607.217
$ws.Range("H141").Value = 3592.6726
$ws.Range("I141").Value = 1603.94
$ws.Range("K141").Value = 4811.82
$ws.Range("M141").Value = 368.1800000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10330.044
$ws.Range("I32").Value = 10574.813
$ws.Range("K32").Value = 10574.813
$ws.Range("M32").Value = -10287.813
$ws.Range("H61").Value = 2462.4062
$ws.Range("I61").Value = 1474.0834
$ws.Range("J61").Value = 5427.375
$ws.Range("K61").Value = 1474.0834
$ws.Range("L61").Value = 5427.375
$ws.Range("M61").Value = -1262.0834
$ws.Range("N61").Value = -5851.375
$ws.Range("H74").Value = 842.96075
$ws.Range("I74").Value = 969.7941
$ws.Range("J74").Value = 589.2941
$ws.Range("K74").Value = 969.7941
$ws.Range("L74").Value = 589.2941
$ws.Range("M74").Value = -95.79409999999996
$ws.Range("N74").Value = -2337.2941
$ws.Range("H77").Value = 842.96075
$ws.Range("I77").Value = 969.7941
$ws.Range("J77").Value = 589.2941
$ws.Range("K77").Value = 4848.970499999999
$ws.Range("L77").Value = 2946.4705
$ws.Range("M77").Value = -480.9704999999994
$ws.Range("N77").Value = -11682.4705
$ws.Range("H132").Value = 1632.0121
$ws.Range("I132").Value = 1302.9166
$ws.Range("J132").Value = 2490.5217
$ws.Range("K132").Value = 3908.7498
$ws.Range("L132").Value = 7471.5651
$ws.Range("M132").Value = -1378.7498
$ws.Range("N132").Value = -12531.5651
$ws.Range("H136").Value = 2462.4062
$ws.Range("I136").Value = 1474.0834
$ws.Range("J136").Value = 5427.375
$ws.Range("K136").Value = 4422.2502
$ws.Range("L136").Value = 16282.125
$ws.Range("M136").Value = -1872.2502
$ws.Range("N136").Value = -21382.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2199.451
$ws.Range("I134").Value = 2026.6562
$ws.Range("J134").Value = 2490.4736
$ws.Range("K134").Value = 6079.9686
$ws.Range("L134").Value = 7471.4208
$ws.Range("M134").Value = -3544.9686
$ws.Range("N134").Value = -12541.4208

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1676.2678
$ws.Range("I31").Value = 1307.3914
$ws.Range("J31").Value = 3373.1
$ws.Range("K31").Value = 1307.3914
$ws.Range("L31").Value = 3373.1
$ws.Range("M31").Value = -1012.3914
$ws.Range("N31").Value = -3963.1
$ws.Range("H34").Value = 1676.2678
$ws.Range("I34").Value = 1307.3914
$ws.Range("J34").Value = 3373.1
$ws.Range("K34").Value = 1307.3914
$ws.Range("L34").Value = 3373.1
$ws.Range("M34").Value = -1105.3914
$ws.Range("N34").Value = -3777.1
$ws.Range("H58").Value = 662640.8
$ws.Range("I58").Value = 950750.7
$ws.Range("J58").Value = 1682.9412
$ws.Range("K58").Value = 950750.7
$ws.Range("L58").Value = 1682.9412
$ws.Range("M58").Value = -950547.7
$ws.Range("N58").Value = -2088.9412
$ws.Range("H132").Value = 246966.2
$ws.Range("I132").Value = 347413.84
$ws.Range("J132").Value = 2125.0625
$ws.Range("K132").Value = 1042241.52
$ws.Range("L132").Value = 6375.1875
$ws.Range("M132").Value = -1039711.52
$ws.Range("N132").Value = -11435.1875
$ws.Range("H134").Value = 1030.0986
$ws.Range("I134").Value = 844.2909
$ws.Range("J134").Value = 1668.8125
$ws.Range("K134").Value = 2532.8727
$ws.Range("L134").Value = 5006.4375
$ws.Range("M134").Value = 2.127300000000105
$ws.Range("N134").Value = -10076.4375
$ws.Range("H136").Value = 662640.8
$ws.Range("I136").Value = 950750.7
$ws.Range("J136").Value = 1682.9412
$ws.Range("K136").Value = 2852252.1
$ws.Range("L136").Value = 5048.8236
$ws.Range("M136").Value = -2849702.1
$ws.Range("N136").Value = -10148.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 6444
$ws.Range("I141").Value = 6444
$ws.Range("K141").Value = 19332
$ws.Range("M141").Value = -14152

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5940
$ws.Range("I70").Value = 5870.7
$ws.Range("J70").Value = 6055.5
$ws.Range("K70").Value = 5870.7
$ws.Range("L70").Value = 6055.5
$ws.Range("M70").Value = -5600.7
$ws.Range("N70").Value = -6595.5
$ws.Range("H73").Value = 5940
$ws.Range("I73").Value = 5870.7
$ws.Range("J73").Value = 6055.5
$ws.Range("K73").Value = 5870.7
$ws.Range("L73").Value = 6055.5
$ws.Range("M73").Value = -4934.7
$ws.Range("N73").Value = -7927.5
$ws.Range("H132").Value = 1657.9333
$ws.Range("I132").Value = 1156.138
$ws.Range("J132").Value = 2567.4375
$ws.Range("K132").Value = 3468.414
$ws.Range("L132").Value = 7702.3125
$ws.Range("M132").Value = -938.4139999999998
$ws.Range("N132").Value = -12762.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 18254.545
$ws.Range("I4").Value = 18633.334
$ws.Range("J4").Value = 18112.5
$ws.Range("K4").Value = 18633.334
$ws.Range("L4").Value = 18112.5
$ws.Range("M4").Value = -18520.334
$ws.Range("N4").Value = -18338.5
$ws.Range("H25").Value = 229500
$ws.Range("I25").Value = 229500
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 229500
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -229270
$ws.Range("N25").ClearContents()
$ws.Range("H28").Value = 18254.545
$ws.Range("I28").Value = 18633.334
$ws.Range("J28").Value = 18112.5
$ws.Range("K28").Value = 18633.334
$ws.Range("L28").Value = 18112.5
$ws.Range("M28").Value = -18401.334
$ws.Range("N28").Value = -18576.5
$ws.Range("H37").Value = 18254.545
$ws.Range("I37").Value = 18633.334
$ws.Range("J37").Value = 18112.5
$ws.Range("K37").Value = 18633.334
$ws.Range("L37").Value = 18112.5
$ws.Range("M37").Value = -18526.334
$ws.Range("N37").Value = -18326.5
$ws.Range("H40").Value = 3668.476
$ws.Range("I40").Value = 3976.182
$ws.Range("J40").Value = 3330
$ws.Range("K40").Value = 3976.182
$ws.Range("L40").Value = 3330
$ws.Range("M40").Value = -3840.182
$ws.Range("N40").Value = -3602
$ws.Range("H132").Value = 1740.4626
$ws.Range("I132").Value = 1603.4426
$ws.Range("J132").Value = 3133.5
$ws.Range("K132").Value = 4810.3278
$ws.Range("L132").Value = 9400.5
$ws.Range("M132").Value = -2280.3278
$ws.Range("N132").Value = -14460.5
$ws.Range("H136").Value = 1509.0695
$ws.Range("I136").Value = 1230.1177
$ws.Range("J136").Value = 6251.25
$ws.Range("K136").Value = 3690.3531
$ws.Range("L136").Value = 18753.75
$ws.Range("M136").Value = -1140.3531
$ws.Range("N136").Value = -23853.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 952.51514
$ws.Range("I132").Value = 640.5769
$ws.Range("J132").Value = 2111.1428
$ws.Range("K132").Value = 1921.7307
$ws.Range("L132").Value = 6333.428400000001
$ws.Range("M132").Value = 608.2692999999999
$ws.Range("N132").Value = -11393.4284
$ws.Range("H136").Value = 900
$ws.Range("I136").Value = 943.1905
$ws.Range("J136").Value = 698.44446
$ws.Range("K136").Value = 2829.5715
$ws.Range("L136").Value = 2095.33338
$ws.Range("M136").Value = -279.5715
$ws.Range("N136").Value = -7195.33338
